$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the result text in B16: "Harus ada aplikasi..." -> "Dihasilkan aplikasi..."
$ws.Range("B16").Value = "Dihasilkan aplikasi yang bisa membantu pencatatan aktivitas keuangan untuk meminimalkan kesalahan pencatatan dan menambah efisiensi waktu."

# Update the view: scroll so A10 is the top-left visible cell, and select B17
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B17").Select() | Out-Null
